$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamps = @(46049,46049.01041666666,46049.02083333334,46049.03125,46049.04166666666,46049.05208333334,46049.0625,46049.07291666666,46049.08333333334,46049.09375,46049.10416666666,46049.11458333334,46049.125,46049.13541666666,46049.14583333334,46049.15625,46049.16666666666,46049.17708333334,46049.1875,46049.19791666666,46049.20833333334,46049.21875,46049.22916666666,46049.23958333334,46049.25,46049.26041666666,46049.27083333334,46049.28125,46049.29166666666,46049.30208333334,46049.3125,46049.32291666666,46049.33333333334,46049.34375,46049.35416666666,46049.36458333334,46049.375,46049.38541666666,46049.39583333334,46049.40625,46049.41666666666,46049.42708333334,46049.4375,46049.44791666666,46049.45833333334,46049.46875,46049.47916666666,46049.48958333334,46049.5,46049.51041666666,46049.52083333334,46049.53125,46049.54166666666,46049.55208333334,46049.5625,46049.57291666666,46049.58333333334,46049.59375,46049.60416666666,46049.61458333334,46049.625,46049.63541666666,46049.64583333334,46049.65625,46049.66666666666,46049.67708333334,46049.6875,46049.69791666666,46049.70833333334,46049.71875,46049.72916666666,46049.73958333334,46049.75,46049.76041666666,46049.77083333334,46049.78125,46049.79166666666,46049.80208333334,46049.8125,46049.82291666666,46049.83333333334,46049.84375,46049.85416666666,46049.86458333334,46049.875,46049.88541666666,46049.89583333334,46049.90625,46049.91666666666,46049.92708333334,46049.9375,46049.94791666666,46049.95833333334,46049.96875,46049.97916666666,46049.98958333334)
$newValues = @(6120,6090,6060,6030,5990,5990,5980,5980,5990,6000,6010,6020,6030,6040,6050,6090,6160,6250,6350,6490,6660,6840,7030,7230,7430,7620,7810,7980,8140,8270,8380,8470,8530,8580,8600,8600,8590,8570,8550,8530,8510,8490,8480,8460,8450,8440,8440,8440,8450,8450,8450,8440,8400,8380,8370,8360,8340,8320,8320,8320,8320,8330,8350,8370,8390,8410,8440,8480,8530,8550,8550,8550,8500,8460,8430,8400,8340,8270,8220,8130,8000,7870,7770,7640,7470,7320,7160,7020,6870,6720,6630,6520,6580,6520,6460,6420)

for ($r = 2; $r -le 97; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $newTimestamps[$r - 2]
    $ws.Cells.Item($r, 2).Value2 = $newValues[$r - 2]
}
